$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Sheet "27.02" (index 5) - add the detailed breakdown data (columns B-I)
# -------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(5)

# Row 2 - Lombardia
$ws.Range("B2").Value = 172
$ws.Range("C2").Value = 41
$ws.Range("D2").Value = 136
$ws.Range("E2").Formula = "=SUM(B2:D2)"
$ws.Range("F2").Value = 40
$ws.Range("G2").Value = 14
$ws.Range("H2").Formula = "=SUM(E2:G2)"
$ws.Range("I2").Value = 3320

# Row 3 - Veneto
$ws.Range("B3").Value = 19
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 82
$ws.Range("E3").Formula = "=SUM(B3:D3)"
$ws.Range("G3").Value = 2
$ws.Range("H3").Formula = "=SUM(E3:G3)"
$ws.Range("I3").Value = 6164

# Row 4 - Emilia Romagna
$ws.Range("B4").Value = 36
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 54
$ws.Range("E4").Formula = "=SUM(B4:D4)"
$ws.Range("G4").Value = 1
$ws.Range("H4").Formula = "=SUM(E4:G4)"
$ws.Range("I4").Value = 1033

# Row 5 - Liguria
$ws.Range("B5").Value = 9
$ws.Range("D5").Value = 10
$ws.Range("E5").Formula = "=SUM(B5:D5)"
$ws.Range("H5").Formula = "=SUM(E5:G5)"
$ws.Range("I5").Value = 78

# Row 6 - Piemonte
$ws.Range("B6").Value = 2
$ws.Range("E6").Formula = "=SUM(B6:D6)"
$ws.Range("H6").Formula = "=SUM(E6:G6)"
$ws.Range("I6").Value = 156

# Row 7 - Toscana
$ws.Range("B7").Value = 2
$ws.Range("E7").Formula = "=SUM(B7:D7)"
$ws.Range("H7").Formula = "=SUM(E7:G7)"
$ws.Range("I7").Value = 410

# Row 8 - Marche
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 1
$ws.Range("E8").Formula = "=SUM(B8:D8)"
$ws.Range("H8").Formula = "=SUM(E8:G8)"
$ws.Range("I8").Value = 46

# Row 9 - Sicilia
$ws.Range("B9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Formula = "=SUM(B9:D9)"
$ws.Range("F9").Value = 2
$ws.Range("H9").Formula = "=SUM(E9:G9)"
$ws.Range("I9").Value = 5

# Row 10 - Lazio
$ws.Range("E10").Formula = "=SUM(B10:D10)"
$ws.Range("F10").Value = 3
$ws.Range("H10").Formula = "=SUM(E10:G10)"
$ws.Range("I10").Value = 552

# Row 11 - Campania
$ws.Range("B11").Value = 2
$ws.Range("D11").Value = 1
$ws.Range("E11").Formula = "=SUM(B11:D11)"
$ws.Range("H11").Formula = "=SUM(E11:G11)"
$ws.Range("I11").Value = 10

# Row 12 - Puglia
$ws.Range("B12").Value = 1
$ws.Range("E12").Formula = "=SUM(B12:D12)"
$ws.Range("H12").Formula = "=SUM(E12:G12)"

# Row 13 - Bolzano
$ws.Range("B13").Value = 1
$ws.Range("E13").Formula = "=SUM(B13:D13)"
$ws.Range("H13").Formula = "=SUM(E13:G13)"
$ws.Range("I13").Value = 2

# Row 14 - Abruzzo
$ws.Range("B14").Value = 1
$ws.Range("E14").Formula = "=SUM(B14:D14)"
$ws.Range("H14").Formula = "=SUM(E14:G14)"
$ws.Range("I14").Value = 33

# Rows 15-22 - only the running total (H) and the tamponi count (I) are known
$ws.Range("H15").Formula = "=SUM(E15:G15)"
$ws.Range("I15").Value = 14

$ws.Range("H16").Formula = "=SUM(E16:G16)"
$ws.Range("I16").Value = 1

$ws.Range("H17").Formula = "=SUM(E17:G17)"
$ws.Range("I17").Value = 8

$ws.Range("H18").Formula = "=SUM(E18:G18)"
$ws.Range("I18").Value = 9

$ws.Range("H19").Formula = "=SUM(E19:G19)"
$ws.Range("I19").Value = 141

$ws.Range("H20").Formula = "=SUM(E20:G20)"
$ws.Range("I20").Value = 32

$ws.Range("H21").Formula = "=SUM(E21:G21)"

$ws.Range("H22").Formula = "=SUM(E22:G22)"

# Totals row (23) - turn the plain F/G totals into SUM formulas and add the
# missing H/I totals, matching the formatting already used by E23
$ws.Range("F23").Formula = "=SUM(F2:F22)"
$ws.Range("G23").Formula = "=SUM(G2:G22)"
$ws.Range("H23").Formula = "=SUM(H2:H22)"
$ws.Range("I23").Formula = "=SUM(I2:I22)"

$ws.Range("E23").Copy()
$ws.Range("F23:I23").PasteSpecial(-4122)

# A handful of cells (D9, D11, B12, B14) carry a distinct font in the
# original edit - recreate that font once and copy its formatting onto the
# other affected cells
$ws.Range("D9").Font.Name = "Arial"
$ws.Range("D9").Font.Size = 12
$ws.Range("D9").Font.ThemeColor = 1
$ws.Range("D9").WrapText = $false

$ws.Range("D9").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B14").PasteSpecial(-4122)

# New column F is now populated - widen it to fit (matches the bestFit width
# Excel computed for the new data)
$ws.Columns.Item(6).ColumnWidth = 11.285714285714286

# -------------------------------------------------------------------------
# Sheet "26.02" (index 4) - selection moved from I23 to B2:B13
# -------------------------------------------------------------------------
$ws26 = $wb.Worksheets.Item(4)
$ws26.Range("B2:B13").Select()

# -------------------------------------------------------------------------
# Make "27.02" the active sheet/tab, with the new selection at I17 - this
# also clears tabSelected/topLeftCell from whichever sheet was active before
# (29.02 in the source workbook)
# -------------------------------------------------------------------------
$ws.Activate()
$ws.Range("I17").Select()
